$wb = $excel.ActiveWorkbook

# --- Sheet "ALL": fix test-id typos (OPQA1366.. -> OPQA-1366..) ---
$wsAll = $wb.Worksheets.Item("ALL")
$wsAll.Range("A3").Value  = "OPQA-1366"
$wsAll.Range("A4").Value  = "OPQA-1367"
$wsAll.Range("A5").Value  = "OPQA-1368"
$wsAll.Range("A6").Value  = "OPQA-1369"
$wsAll.Range("A7").Value  = "OPQA-1370"
$wsAll.Range("A8").Value  = "OPQA-1371"
$wsAll.Range("A9").Value  = "OPQA-1372"
$wsAll.Range("A10").Value = "OPQA-1373"
$wsAll.Range("A12").Value = "OPQA-1374"
$wsAll.Range("A13").Value = "OPQA-1375"

# --- Sheet "People" ---
$wsPeople = $wb.Worksheets.Item("People")
$wsPeople.Range("A2").Value = "OPQA-1376"
$wsPeople.Range("A3").Value = "OPQA-1377"
$wsPeople.Range("A4").Value = "OPQA-1378"

# --- Sheet "Patents" ---
$wsPatents = $wb.Worksheets.Item("Patents")
$wsPatents.Range("A2").Value = "OPQA-1379"
$wsPatents.Range("A4").Value = "OPQA-1380"

# --- Sheet "Articles" ---
$wsArticles = $wb.Worksheets.Item("Articles")
$wsArticles.Range("A2").Value = "OPQA-1381"
$wsArticles.Range("A4").Value = "OPQA-1382"

# --- Sheet "Posts" ---
$wsPosts = $wb.Worksheets.Item("Posts")
$wsPosts.Range("A2").Value = "OPQA-1383"
$wsPosts.Range("A4").Value = "OPQA-1384"

# --- Sheet "Profile" ---
$wsProfile = $wb.Worksheets.Item("Profile")
$wsProfile.Range("A2").Value = "OPQA-1384"
$wsProfile.Range("A3").Value = "OPQA-1385"

# --- Update each sheet's view / selection (also clears any stale topLeftCell) ---
$wsAll.Activate()
$wsAll.Range("A12:A13").Select()

$wsPeople.Activate()
$wsPeople.Range("A2:A4").Select()

$wsPatents.Activate()
$wsPatents.Range("A4").Select()

$wsArticles.Activate()
$wsArticles.Range("A4").Select()

$wsPosts.Activate()
$wsPosts.Range("A4").Select()

# Profile ends up the active / selected tab, matching activeTab="5" in the workbook views
$wsProfile.Activate()
$wsProfile.Range("A2:A3").Select()
